# Bubble-sort figure: the tracked element "25" becomes "4" as it is
# swapped/bubbled through the array, together with the matching
# "Swap 109 and 25;" annotation between the 2nd and 3rd table.
#
# Slide 1 layout (document order):
#   Shape 1  -> Table (row: 109, -13, 25, -48, 25)   -- col 3 changes
#   Shape 8  -> Table (row: -13, 109, 25, -48, 25)   -- col 3 changes
#   Shape 12 -> TextBox "Constraint not satisfied: / Swap 109 and 25; / Increment Index"
#   Shape 13 -> Table (row: -13, 25, 109, -48, 25)   -- col 2 changes
#   Shape 18 -> Table (row: -13, 25, -48, 109, 25)   -- col 2 changes
#   Shape 23 -> Table (row: -13, 25, -48, 25, 109)   -- col 2 changes (NOT touched)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Table 1 (shape 1), column 3: 25 -> 4
$s.Shapes.Item(1).Table.Cell(1, 3).Shape.TextFrame.TextRange.Text = "4"

# Table 2 (shape 8), column 3: 25 -> 4
$s.Shapes.Item(8).Table.Cell(1, 3).Shape.TextFrame.TextRange.Text = "4"

# TextBox (shape 12): update the "Swap 109 and 25;" paragraph only,
# replacing the whole paragraph range so the run stays intact.
$tb = $s.Shapes.Item(12)
$tr = $tb.TextFrame.TextRange
$full = $tr.Text
$firstCr = $full.IndexOf([char]13)
$secondCr = $full.IndexOf([char]13, $firstCr + 1)
$paraStart = $firstCr + 2
$paraLen = $secondCr - $firstCr - 1
$tr.Characters($paraStart, $paraLen).Text = "Swap 109 and 4; "

# Table 3 (shape 13), column 2: 25 -> 4
$s.Shapes.Item(13).Table.Cell(1, 2).Shape.TextFrame.TextRange.Text = "4"

# Table 4 (shape 18), column 2: 25 -> 4
$s.Shapes.Item(18).Table.Cell(1, 2).Shape.TextFrame.TextRange.Text = "4"

# Table 5 (shape 23), column 2: 25 -> 4
$s.Shapes.Item(23).Table.Cell(1, 2).Shape.TextFrame.TextRange.Text = "4"
